$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 5381.1113
$ws.Range("I74").Value = 5718.5713
$ws.Range("K74").Value = 5718.5713
$ws.Range("M74").Value = -4782.5713

$ws.Range("H77").Value = 5381.1113
$ws.Range("I77").Value = 5718.5713
$ws.Range("K77").Value = 28592.8565
$ws.Range("M77").Value = -23912.8565

$ws.Range("H81").Value = 35000
$ws.Range("J81").Value = 35000
$ws.Range("L81").Value = 35000
$ws.Range("N81").Value = -36996

$ws.Range("H84").Value = 35000
$ws.Range("J84").Value = 35000
$ws.Range("L84").Value = 105000
$ws.Range("N84").Value = -114984

$ws.Range("H137").Value = 1560.4117
$ws.Range("I137").Value = 1220.3077
$ws.Range("J137").Value = 2665.75
$ws.Range("K137").Value = 3660.9231
$ws.Range("L137").Value = 7997.25
$ws.Range("M137").Value = -1110.9231
$ws.Range("N137").Value = -13097.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4417.6597
$ws.Range("I32").Value = 2442.3428
$ws.Range("J32").Value = 10179
$ws.Range("K32").Value = 2442.3428
$ws.Range("L32").Value = 10179
$ws.Range("M32").Value = -2155.3428
$ws.Range("N32").Value = -10753

$ws.Range("H61").Value = 4095.4167
$ws.Range("I61").Value = 2422.2222
$ws.Range("K61").Value = 2422.2222
$ws.Range("M61").Value = -2210.2222

$ws.Range("H74").Value = 693.36365
$ws.Range("I74").Value = 693.36365
$ws.Range("K74").Value = 693.36365
$ws.Range("M74").Value = 180.63635

$ws.Range("H77").Value = 693.36365
$ws.Range("I77").Value = 693.36365
$ws.Range("K77").Value = 3466.81825
$ws.Range("M77").Value = 901.1817499999997

$ws.Range("H122").Value = 2250.923
$ws.Range("I122").Value = 2162.4443
$ws.Range("J122").Value = 2450
$ws.Range("K122").Value = 6487.3329
$ws.Range("L122").Value = 7350
$ws.Range("M122").Value = -4037.3329
$ws.Range("N122").Value = -12250

$ws.Range("H136").Value = 4095.4167
$ws.Range("I136").Value = 2422.2222
$ws.Range("K136").Value = 7266.6666
$ws.Range("M136").Value = -4716.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 424.8
$ws.Range("I22").Value = 387.25
$ws.Range("J22").Value = 575
$ws.Range("K22").Value = 387.25
$ws.Range("L22").Value = 575
$ws.Range("M22").Value = -214.25
$ws.Range("N22").Value = -921

$ws.Range("H62").Value = 10000
$ws.Range("I62").Value = 10000
$ws.Range("K62").Value = 10000
$ws.Range("M62").Value = -9314

$ws.Range("H65").Value = 10000
$ws.Range("I65").Value = 10000
$ws.Range("K65").Value = 30000
$ws.Range("M65").Value = -26568

$ws.Range("H94").Value = 759.2857
$ws.Range("I94").Value = 739.7273
$ws.Range("K94").Value = 739.7273
$ws.Range("M94").Value = -288.7273

$ws.Range("H100").Value = 33000
$ws.Range("J100").Value = 33000
$ws.Range("L100").Value = 33000
$ws.Range("N100").Value = -35164

$ws.Range("H126").Value = 38421.05
$ws.Range("J126").Value = 38421.05
$ws.Range("L126").Value = 38421.05
$ws.Range("N126").Value = -48301.05

$ws.Range("H134").Value = 12960.296
$ws.Range("I134").Value = 13414.869
$ws.Range("K134").Value = 40244.607
$ws.Range("M134").Value = -37709.607

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3203.45
$ws.Range("I31").Value = 1531.7646
$ws.Range("J31").Value = 4439.0435
$ws.Range("K31").Value = 1531.7646
$ws.Range("L31").Value = 4439.0435
$ws.Range("M31").Value = -1236.7646
$ws.Range("N31").Value = -5029.0435

$ws.Range("H34").Value = 3203.45
$ws.Range("I34").Value = 1531.7646
$ws.Range("J34").Value = 4439.0435
$ws.Range("K34").Value = 1531.7646
$ws.Range("L34").Value = 4439.0435
$ws.Range("M34").Value = -1329.7646
$ws.Range("N34").Value = -4843.0435

$ws.Range("H51").Value = 30933.334
$ws.Range("J51").Value = 30933.334
$ws.Range("L51").Value = 30933.334
$ws.Range("N51").Value = -32405.334

$ws.Range("H58").Value = 2175422.5
$ws.Range("I58").Value = 4832147.5
$ws.Range("J58").Value = 1738.3636
$ws.Range("K58").Value = 4832147.5
$ws.Range("L58").Value = 1738.3636
$ws.Range("M58").Value = -4831944.5
$ws.Range("N58").Value = -2144.3636

$ws.Range("H61").Value = 30933.334
$ws.Range("J61").Value = 30933.334
$ws.Range("L61").Value = 30933.334
$ws.Range("N61").Value = -31629.334

$ws.Range("H94").Value = 933
$ws.Range("I94").Value = 756
$ws.Range("J94").Value = 977.25
$ws.Range("K94").Value = 756
$ws.Range("L94").Value = 977.25
$ws.Range("M94").Value = -305
$ws.Range("N94").Value = -1879.25

$ws.Range("H134").Value = 2026.875
$ws.Range("I134").Value = 1495.4
$ws.Range("K134").Value = 4486.200000000001
$ws.Range("M134").Value = -1951.200000000001

$ws.Range("H136").Value = 2175422.5
$ws.Range("I136").Value = 4832147.5
$ws.Range("J136").Value = 1738.3636
$ws.Range("K136").Value = 14496442.5
$ws.Range("L136").Value = 5215.0908
$ws.Range("M136").Value = -14493892.5
$ws.Range("N136").Value = -10315.0908

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 10909.3
$ws.Range("I113").Value = 33868
$ws.Range("J113").Value = 1069.8572
$ws.Range("K113").Value = 101604
$ws.Range("L113").Value = 3209.5716
$ws.Range("M113").Value = -99434
$ws.Range("N113").Value = -7549.571599999999

$ws.Range("H131").Value = 9093.833000000001
$ws.Range("I131").Value = 463.25
$ws.Range("J131").Value = 10002.315
$ws.Range("K131").Value = 1389.75
$ws.Range("L131").Value = 30006.945
$ws.Range("M131").Value = 3650.25
$ws.Range("N131").Value = -40086.945

$ws.Range("H133").Value = 3155

$ws.Range("H134").Value = 1633
$ws.Range("I134").Value = 1259.579
$ws.Range("K134").Value = 3778.737
$ws.Range("M134").Value = 1291.263

$ws.Range("H136").Value = 3143
$ws.Range("I136").Value = 3165
$ws.Range("J136").Value = 3033
$ws.Range("K136").Value = 9495
$ws.Range("L136").Value = 9099
$ws.Range("M136").Value = -4395
$ws.Range("N136").Value = -19299

$ws.Range("H139").Value = 17094
$ws.Range("I139").Value = 19610
$ws.Range("K139").Value = 58830
$ws.Range("M139").Value = -53690

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2307
$ws.Range("I80").Value = 1930
$ws.Range("K80").Value = 1930
$ws.Range("M80").Value = -932

$ws.Range("H83").Value = 2307
$ws.Range("I83").Value = 1930
$ws.Range("K83").Value = 9650
$ws.Range("M83").Value = -4658

$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

$ws.Range("H132").Value = 1752054
$ws.Range("I132").Value = 2566156.5
$ws.Range("K132").Value = 7698469.5
$ws.Range("M132").Value = -7695939.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1842.1786
$ws.Range("I132").Value = 1718.8
$ws.Range("J132").Value = 1869
$ws.Range("K132").Value = 5156.4
$ws.Range("L132").Value = 5607
$ws.Range("M132").Value = -2626.4
$ws.Range("N132").Value = -10667

$ws.Range("H136").Value = 4312.1875
$ws.Range("I136").Value = 2185.8572
$ws.Range("J136").Value = 5966
$ws.Range("K136").Value = 6557.571599999999
$ws.Range("L136").Value = 17898
$ws.Range("M136").Value = -4007.571599999999
$ws.Range("N136").Value = -22998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 11415.667
$ws.Range("J41").Value = 11415.667
$ws.Range("L41").Value = 11415.667
$ws.Range("N41").Value = -12195.667

$ws.Range("H132").Value = 6992.3
$ws.Range("I132").Value = 1471.5
$ws.Range("J132").Value = 8999.862999999999
$ws.Range("K132").Value = 4414.5
$ws.Range("L132").Value = 26999.589
$ws.Range("M132").Value = -1884.5
$ws.Range("N132").Value = -32059.589

$ws.Range("H136").Value = 42738636
$ws.Range("I136").Value = 69448136
$ws.Range("J136").Value = 3440
$ws.Range("K136").Value = 208344408
$ws.Range("L136").Value = 10320
$ws.Range("M136").Value = -208341858
$ws.Range("N136").Value = -15420
